$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the date on the existing last entry (row 39): 01.11.2018 -> 02.11.2018
$ws.Range("F39").Value = 43406

# New entry: "Architektur umgebaut" / "5h" / 02.11.2018
$ws.Range("B41").Value = "Architektur umgebaut"
$ws.Range("D41").Value = "5h"
$ws.Range("F41").Value = 43406
$ws.Range("F41").NumberFormat = "d-mmm"

# New entry: "Projektilsystem funktionsfähig" / "2h" / 02.11.2018
$ws.Range("B42").Value = "Projektilsystem funktionsfähig"
$ws.Range("D42").Value = "2h"
$ws.Range("F42").Value = 43406
$ws.Range("F42").NumberFormat = "d-mmm"

# Update the selected cell to reflect the newly added data area
$ws.Range("F13").Select()
